$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("Overall")
$ws.Range("A1").Value = "Share of 990 filers with government grants at risk"
$ws.Range("B1").Value = "Number of 990 filers with government grants"
$ws.Range("C1").Value = "Total government grants (`$)"
$ws.Range("D1").Value = "Size of operating surplus with government grants"
$ws.Range("E1").Value = "Size of operating surplus without government grants"
$ws.Range("A2").Value = "'60.30%"
$ws.Range("A2").Style = "Normal"
$ws.Range("B2").Value = "'3,018"
$ws.Range("B2").Style = "Normal"
$ws.Range("C2").Value = "'`$10,165,747,853"
$ws.Range("C2").Style = "Normal"
$ws.Range("D2").Value = "'11.74%"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "'-5.38%"
$ws.Range("E2").Style = "Normal"

$ws = $wb.Worksheets.Item("County")
$ws.Range("A1").Value = "Geography"
$ws.Range("B1").Value = "Share of 990 filers with government grants at risk"
$ws.Range("C1").Value = "Number of 990 filers with government grants"
$ws.Range("D1").Value = "Total government grants (`$)"
$ws.Range("E1").Value = "Size of operating surplus with government grants"
$ws.Range("F1").Value = "Size of operating surplus without government grants"
$ws.Range("A2").Value = "United States"
$ws.Range("B2").Value = "'67.35%"
$ws.Range("B2").Style = "Normal"
$ws.Range("C2").Value = "'103,475"
$ws.Range("C2").Style = "Normal"
$ws.Range("D2").Value = "'`$267,700,640,005"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "'9.05%"
$ws.Range("E2").Style = "Normal"
$ws.Range("F2").Value = "'-12.83%"
$ws.Range("F2").Style = "Normal"
$ws.Range("A3").Value = "Massachusetts"
$ws.Range("B3").Value = "'60.30%"
$ws.Range("B3").Style = "Normal"
$ws.Range("C3").Value = "'3,018"
$ws.Range("C3").Style = "Normal"
$ws.Range("D3").Value = "'`$10,165,747,853"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "'11.74%"
$ws.Range("E3").Style = "Normal"
$ws.Range("F3").Value = "'-5.38%"
$ws.Range("F3").Style = "Normal"
$ws.Range("A4").Value = "Barnstable County"
$ws.Range("B4").Value = "'65.56%"
$ws.Range("B4").Style = "Normal"
$ws.Range("C4").Value = "'151"
$ws.Range("C4").Style = "Normal"
$ws.Range("D4").Value = "'`$423,474,552"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "'12.61%"
$ws.Range("E4").Style = "Normal"
$ws.Range("F4").Value = "'-10.35%"
$ws.Range("F4").Style = "Normal"
$ws.Range("A5").Value = "Berkshire County"
$ws.Range("B5").Value = "'49.58%"
$ws.Range("B5").Style = "Normal"
$ws.Range("C5").Value = "'119"
$ws.Range("C5").Style = "Normal"
$ws.Range("D5").Value = "'`$111,463,153"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "'15.33%"
$ws.Range("E5").Style = "Normal"
$ws.Range("F5").Value = "'0.77%"
$ws.Range("F5").Style = "Normal"
$ws.Range("A6").Value = "Bristol County"
$ws.Range("B6").Value = "'66.67%"
$ws.Range("B6").Style = "Normal"
$ws.Range("C6").Value = "'147"
$ws.Range("C6").Style = "Normal"
$ws.Range("D6").Value = "'`$367,318,906"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "'10.77%"
$ws.Range("E6").Style = "Normal"
$ws.Range("F6").Value = "'-8.84%"
$ws.Range("F6").Style = "Normal"
$ws.Range("A7").Value = "Dukes County"
$ws.Range("B7").Value = "'46.67%"
$ws.Range("B7").Style = "Normal"
$ws.Range("C7").Value = "'30"
$ws.Range("C7").Style = "Normal"
$ws.Range("D7").Value = "'`$12,840,897"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "'24.19%"
$ws.Range("E7").Style = "Normal"
$ws.Range("F7").Value = "'5.07%"
$ws.Range("F7").Style = "Normal"
$ws.Range("A8").Value = "Essex County"
$ws.Range("B8").Value = "'62.36%"
$ws.Range("B8").Style = "Normal"
$ws.Range("C8").Value = "'271"
$ws.Range("C8").Style = "Normal"
$ws.Range("D8").Value = "'`$573,711,863"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "'11.25%"
$ws.Range("E8").Style = "Normal"
$ws.Range("F8").Value = "'-6.04%"
$ws.Range("F8").Style = "Normal"
$ws.Range("A9").Value = "Franklin County"
$ws.Range("B9").Value = "'74.29%"
$ws.Range("B9").Style = "Normal"
$ws.Range("C9").Value = "'35"
$ws.Range("C9").Style = "Normal"
$ws.Range("D9").Value = "'`$34,705,521"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "'7.03%"
$ws.Range("E9").Style = "Normal"
$ws.Range("F9").Value = "'-25.52%"
$ws.Range("F9").Style = "Normal"
$ws.Range("A10").Value = "Hampden County"
$ws.Range("B10").Value = "'64.71%"
$ws.Range("B10").Style = "Normal"
$ws.Range("C10").Value = "'153"
$ws.Range("C10").Style = "Normal"
$ws.Range("D10").Value = "'`$642,620,761"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "'11.47%"
$ws.Range("E10").Style = "Normal"
$ws.Range("F10").Value = "'-12.38%"
$ws.Range("F10").Style = "Normal"
$ws.Range("A11").Value = "Hampshire County"
$ws.Range("B11").Value = "'62.50%"
$ws.Range("B11").Style = "Normal"
$ws.Range("C11").Value = "'96"
$ws.Range("C11").Style = "Normal"
$ws.Range("D11").Value = "'`$184,336,565"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "'13.59%"
$ws.Range("E11").Style = "Normal"
$ws.Range("F11").Value = "'-8.33%"
$ws.Range("F11").Style = "Normal"
$ws.Range("A12").Value = "Middlesex County"
$ws.Range("B12").Value = "'55.41%"
$ws.Range("B12").Style = "Normal"
$ws.Range("C12").Value = "'702"
$ws.Range("C12").Style = "Normal"
$ws.Range("D12").Value = "'`$3,522,554,953"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "'13.29%"
$ws.Range("E12").Style = "Normal"
$ws.Range("F12").Value = "'-2.31%"
$ws.Range("F12").Style = "Normal"
$ws.Range("A13").Value = "Nantucket County"
$ws.Range("B13").Value = "'38.71%"
$ws.Range("B13").Style = "Normal"
$ws.Range("C13").Value = "'31"
$ws.Range("C13").Style = "Normal"
$ws.Range("D13").Value = "'`$7,595,704"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "'31.06%"
$ws.Range("E13").Style = "Normal"
$ws.Range("F13").Value = "'6.40%"
$ws.Range("F13").Style = "Normal"
$ws.Range("A14").Value = "Norfolk County"
$ws.Range("B14").Value = "'59.01%"
$ws.Range("B14").Style = "Normal"
$ws.Range("C14").Value = "'222"
$ws.Range("C14").Style = "Normal"
$ws.Range("D14").Value = "'`$573,629,276"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "'11.38%"
$ws.Range("E14").Style = "Normal"
$ws.Range("F14").Value = "'-2.72%"
$ws.Range("F14").Style = "Normal"
$ws.Range("A15").Value = "Plymouth County"
$ws.Range("B15").Value = "'55.08%"
$ws.Range("B15").Style = "Normal"
$ws.Range("C15").Value = "'118"
$ws.Range("C15").Style = "Normal"
$ws.Range("D15").Value = "'`$198,197,463"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "'11.61%"
$ws.Range("E15").Style = "Normal"
$ws.Range("F15").Value = "'-2.82%"
$ws.Range("F15").Style = "Normal"
$ws.Range("A16").Value = "Suffolk County"
$ws.Range("B16").Value = "'63.23%"
$ws.Range("B16").Style = "Normal"
$ws.Range("C16").Value = "'718"
$ws.Range("C16").Style = "Normal"
$ws.Range("D16").Value = "'`$2,909,376,536"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "'10.84%"
$ws.Range("E16").Style = "Normal"
$ws.Range("F16").Value = "'-7.96%"
$ws.Range("F16").Style = "Normal"
$ws.Range("A17").Value = "Worcester County"
$ws.Range("B17").Value = "'64.44%"
$ws.Range("B17").Style = "Normal"
$ws.Range("C17").Value = "'225"
$ws.Range("C17").Style = "Normal"
$ws.Range("D17").Value = "'`$603,921,703"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "'8.82%"
$ws.Range("E17").Style = "Normal"
$ws.Range("F17").Value = "'-6.77%"
$ws.Range("F17").Style = "Normal"

$ws = $wb.Worksheets.Item("Congressional District")
$ws.Range("A1").Value = "Geography"
$ws.Range("B1").Value = "Share of 990 filers with government grants at risk"
$ws.Range("C1").Value = "Number of 990 filers with government grants"
$ws.Range("D1").Value = "Total government grants (`$)"
$ws.Range("E1").Value = "Size of operating surplus with government grants"
$ws.Range("F1").Value = "Size of operating surplus without government grants"
$ws.Range("A2").Value = "United States"
$ws.Range("B2").Value = "'67.35%"
$ws.Range("B2").Style = "Normal"
$ws.Range("C2").Value = "'103,475"
$ws.Range("C2").Style = "Normal"
$ws.Range("D2").Value = "'`$267,700,640,005"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "'9.05%"
$ws.Range("E2").Style = "Normal"
$ws.Range("F2").Value = "'-12.83%"
$ws.Range("F2").Style = "Normal"
$ws.Range("A3").Value = "Massachusetts"
$ws.Range("B3").Value = "'60.30%"
$ws.Range("B3").Style = "Normal"
$ws.Range("C3").Value = "'3,018"
$ws.Range("C3").Style = "Normal"
$ws.Range("D3").Value = "'`$10,165,747,853"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "'11.74%"
$ws.Range("E3").Style = "Normal"
$ws.Range("F3").Value = "'-5.38%"
$ws.Range("F3").Style = "Normal"
$ws.Range("A4").Value = "Congressional District 1"
$ws.Range("B4").Value = "'60.52%"
$ws.Range("B4").Style = "Normal"
$ws.Range("C4").Value = "'309"
$ws.Range("C4").Style = "Normal"
$ws.Range("D4").Value = "'`$810,434,787"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "'11.52%"
$ws.Range("E4").Style = "Normal"
$ws.Range("F4").Value = "'-9.00%"
$ws.Range("F4").Style = "Normal"
$ws.Range("A5").Value = "Congressional District 2"
$ws.Range("B5").Value = "'63.64%"
$ws.Range("B5").Style = "Normal"
$ws.Range("C5").Value = "'286"
$ws.Range("C5").Style = "Normal"
$ws.Range("D5").Value = "'`$698,888,291"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "'10.07%"
$ws.Range("E5").Style = "Normal"
$ws.Range("F5").Value = "'-7.08%"
$ws.Range("F5").Style = "Normal"
$ws.Range("A6").Value = "Congressional District 3"
$ws.Range("B6").Value = "'57.83%"
$ws.Range("B6").Style = "Normal"
$ws.Range("C6").Value = "'230"
$ws.Range("C6").Style = "Normal"
$ws.Range("D6").Value = "'`$608,155,140"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "'13.27%"
$ws.Range("E6").Style = "Normal"
$ws.Range("F6").Value = "'-5.24%"
$ws.Range("F6").Style = "Normal"
$ws.Range("A7").Value = "Congressional District 4"
$ws.Range("B7").Value = "'55.43%"
$ws.Range("B7").Style = "Normal"
$ws.Range("C7").Value = "'267"
$ws.Range("C7").Style = "Normal"
$ws.Range("D7").Value = "'`$691,593,479"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "'12.61%"
$ws.Range("E7").Style = "Normal"
$ws.Range("F7").Value = "'-2.14%"
$ws.Range("F7").Style = "Normal"
$ws.Range("A8").Value = "Congressional District 5"
$ws.Range("B8").Value = "'58.82%"
$ws.Range("B8").Style = "Normal"
$ws.Range("C8").Value = "'357"
$ws.Range("C8").Style = "Normal"
$ws.Range("D8").Value = "'`$1,734,751,035"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "'11.79%"
$ws.Range("E8").Style = "Normal"
$ws.Range("F8").Value = "'-3.56%"
$ws.Range("F8").Style = "Normal"
$ws.Range("A9").Value = "Congressional District 6"
$ws.Range("B9").Value = "'62.03%"
$ws.Range("B9").Style = "Normal"
$ws.Range("C9").Value = "'237"
$ws.Range("C9").Style = "Normal"
$ws.Range("D9").Value = "'`$920,998,514"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "'10.95%"
$ws.Range("E9").Style = "Normal"
$ws.Range("F9").Value = "'-5.97%"
$ws.Range("F9").Style = "Normal"
$ws.Range("A10").Value = "Congressional District 7"
$ws.Range("B10").Value = "'62.71%"
$ws.Range("B10").Style = "Normal"
$ws.Range("C10").Value = "'539"
$ws.Range("C10").Style = "Normal"
$ws.Range("D10").Value = "'`$3,069,340,319"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "'11.39%"
$ws.Range("E10").Style = "Normal"
$ws.Range("F10").Value = "'-6.90%"
$ws.Range("F10").Style = "Normal"
$ws.Range("A11").Value = "Congressional District 8"
$ws.Range("B11").Value = "'61.29%"
$ws.Range("B11").Style = "Normal"
$ws.Range("C11").Value = "'434"
$ws.Range("C11").Style = "Normal"
$ws.Range("D11").Value = "'`$1,016,216,002"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "'10.50%"
$ws.Range("E11").Style = "Normal"
$ws.Range("F11").Value = "'-6.30%"
$ws.Range("F11").Style = "Normal"
$ws.Range("A12").Value = "Congressional District 9"
$ws.Range("B12").Value = "'58.22%"
$ws.Range("B12").Style = "Normal"
$ws.Range("C12").Value = "'359"
$ws.Range("C12").Style = "Normal"
$ws.Range("D12").Value = "'`$615,370,286"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "'14.81%"
$ws.Range("E12").Style = "Normal"
$ws.Range("F12").Value = "'-3.67%"
$ws.Range("F12").Style = "Normal"

$ws = $wb.Worksheets.Item("Size")
$ws.Range("A1").Value = "Size"
$ws.Range("B1").Value = "Share of 990 filers with government grants at risk"
$ws.Range("C1").Value = "Number of 990 filers with government grants"
$ws.Range("D1").Value = "Total government grants (`$)"
$ws.Range("E1").Value = "Size of operating surplus with government grants"
$ws.Range("F1").Value = "Size of operating surplus without government grants"
$ws.Range("A2").Value = "Between `$100K and `$499K"
$ws.Range("B2").Value = "'57.77%"
$ws.Range("B2").Style = "Normal"
$ws.Range("C2").Value = "'772"
$ws.Range("C2").Style = "Normal"
$ws.Range("D2").Value = "'`$72,671,079"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "'14.57%"
$ws.Range("E2").Style = "Normal"
$ws.Range("F2").Value = "'-4.66%"
$ws.Range("F2").Style = "Normal"
$ws.Range("A3").Value = "Between `$1M and `$4.99M"
$ws.Range("B3").Value = "'59.81%"
$ws.Range("B3").Style = "Normal"
$ws.Range("C3").Value = "'928"
$ws.Range("C3").Style = "Normal"
$ws.Range("D3").Value = "'`$672,947,375"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "'12.39%"
$ws.Range("E3").Style = "Normal"
$ws.Range("F3").Value = "'-5.88%"
$ws.Range("F3").Style = "Normal"
$ws.Range("A4").Value = "Between `$500K and `$999K"
$ws.Range("B4").Value = "'58.28%"
$ws.Range("B4").Style = "Normal"
$ws.Range("C4").Value = "'489"
$ws.Range("C4").Style = "Normal"
$ws.Range("D4").Value = "'`$111,110,741"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "'13.77%"
$ws.Range("E4").Style = "Normal"
$ws.Range("F4").Value = "'-3.66%"
$ws.Range("F4").Style = "Normal"
$ws.Range("A5").Value = "Between `$5M and `$9.99M"
$ws.Range("B5").Value = "'60.87%"
$ws.Range("B5").Style = "Normal"
$ws.Range("C5").Value = "'230"
$ws.Range("C5").Style = "Normal"
$ws.Range("D5").Value = "'`$432,831,622"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "'8.98%"
$ws.Range("E5").Style = "Normal"
$ws.Range("F5").Value = "'-5.30%"
$ws.Range("F5").Style = "Normal"
$ws.Range("A6").Value = "Greater than `$10M"
$ws.Range("B6").Value = "'68.32%"
$ws.Range("B6").Style = "Normal"
$ws.Range("C6").Value = "'483"
$ws.Range("C6").Style = "Normal"
$ws.Range("D6").Value = "'`$8,871,097,129"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "'6.41%"
$ws.Range("E6").Style = "Normal"
$ws.Range("F6").Value = "'-6.83%"
$ws.Range("F6").Style = "Normal"
$ws.Range("A7").Value = "Less than `$100K"
$ws.Range("B7").Value = "'55.17%"
$ws.Range("B7").Style = "Normal"
$ws.Range("C7").Value = "'116"
$ws.Range("C7").Style = "Normal"
$ws.Range("D7").Value = "'`$5,089,907"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "'28.37%"
$ws.Range("E7").Style = "Normal"
$ws.Range("F7").Value = "'-5.86%"
$ws.Range("F7").Style = "Normal"
$ws.Range("A8").Value = "Total"
$ws.Range("B8").Value = "'60.30%"
$ws.Range("B8").Style = "Normal"
$ws.Range("C8").Value = "'3,018"
$ws.Range("C8").Style = "Normal"
$ws.Range("D8").Value = "'`$10,165,747,853"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "'11.74%"
$ws.Range("E8").Style = "Normal"
$ws.Range("F8").Value = "'-5.38%"
$ws.Range("F8").Style = "Normal"

$ws = $wb.Worksheets.Item("Subsector")
$ws.Range("A1").Value = "Subsector"
$ws.Range("B1").Value = "Share of 990 filers with government grants at risk"
$ws.Range("C1").Value = "Number of 990 filers with government grants"
$ws.Range("D1").Value = "Total government grants (`$)"
$ws.Range("E1").Value = "Size of operating surplus with government grants"
$ws.Range("F1").Value = "Size of operating surplus without government grants"
$ws.Range("A2").Value = "Arts, Culture, and Humanities"
$ws.Range("B2").Value = "'58.17%"
$ws.Range("B2").Style = "Normal"
$ws.Range("C2").Value = "'361"
$ws.Range("C2").Style = "Normal"
$ws.Range("D2").Value = "'`$175,361,193"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "'15.99%"
$ws.Range("E2").Style = "Normal"
$ws.Range("F2").Value = "'-2.53%"
$ws.Range("F2").Style = "Normal"
$ws.Range("A3").Value = "Education (Excluding Universities)"
$ws.Range("B3").Value = "'52.69%"
$ws.Range("B3").Style = "Normal"
$ws.Range("C3").Value = "'372"
$ws.Range("C3").Style = "Normal"
$ws.Range("D3").Value = "'`$599,606,169"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "'15.07%"
$ws.Range("E3").Style = "Normal"
$ws.Range("F3").Value = "'-1.08%"
$ws.Range("F3").Style = "Normal"
$ws.Range("A4").Value = "Environment and Animals"
$ws.Range("B4").Value = "'39.20%"
$ws.Range("B4").Style = "Normal"
$ws.Range("C4").Value = "'125"
$ws.Range("C4").Style = "Normal"
$ws.Range("D4").Value = "'`$86,384,192"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "'25.48%"
$ws.Range("E4").Style = "Normal"
$ws.Range("F4").Value = "'6.33%"
$ws.Range("F4").Style = "Normal"
$ws.Range("A5").Value = "Health (Excluding Hospitals)"
$ws.Range("B5").Value = "'66.67%"
$ws.Range("B5").Style = "Normal"
$ws.Range("C5").Value = "'234"
$ws.Range("C5").Style = "Normal"
$ws.Range("D5").Value = "'`$1,379,136,490"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "'8.98%"
$ws.Range("E5").Style = "Normal"
$ws.Range("F5").Value = "'-8.00%"
$ws.Range("F5").Style = "Normal"
$ws.Range("A6").Value = "Hospitals"
$ws.Range("B6").Value = "'87.50%"
$ws.Range("B6").Style = "Normal"
$ws.Range("C6").Value = "'8"
$ws.Range("C6").Style = "Normal"
$ws.Range("D6").Value = "'`$203,091,046"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "'-0.03%"
$ws.Range("E6").Style = "Normal"
$ws.Range("F6").Value = "'-25.73%"
$ws.Range("F6").Style = "Normal"
$ws.Range("A7").Value = "Human Services"
$ws.Range("B7").Value = "'64.44%"
$ws.Range("B7").Style = "Normal"
$ws.Range("C7").Value = "'852"
$ws.Range("C7").Style = "Normal"
$ws.Range("D7").Value = "'`$2,686,654,793"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "'10.04%"
$ws.Range("E7").Style = "Normal"
$ws.Range("F7").Value = "'-8.83%"
$ws.Range("F7").Style = "Normal"
$ws.Range("A8").Value = "International, Foreign Affairs"
$ws.Range("B8").Value = "'57.78%"
$ws.Range("B8").Style = "Normal"
$ws.Range("C8").Value = "'45"
$ws.Range("C8").Style = "Normal"
$ws.Range("D8").Value = "'`$583,969,709"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "'5.26%"
$ws.Range("E8").Style = "Normal"
$ws.Range("F8").Value = "'-3.09%"
$ws.Range("F8").Style = "Normal"
$ws.Range("A9").Value = "Mutual/Membership Benefit"
$ws.Range("B9").Value = "'100.00%"
$ws.Range("B9").Style = "Normal"
$ws.Range("C9").Value = "'1"
$ws.Range("C9").Style = "Normal"
$ws.Range("D9").Value = "'`$23,120"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "'-90.32%"
$ws.Range("E9").Style = "Normal"
$ws.Range("F9").Value = "'-95.95%"
$ws.Range("F9").Style = "Normal"
$ws.Range("A10").Value = "Public, Societal Benefit"
$ws.Range("B10").Value = "'61.98%"
$ws.Range("B10").Style = "Normal"
$ws.Range("C10").Value = "'242"
$ws.Range("C10").Style = "Normal"
$ws.Range("D10").Value = "'`$419,647,334"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "'11.08%"
$ws.Range("E10").Style = "Normal"
$ws.Range("F10").Value = "'-11.38%"
$ws.Range("F10").Style = "Normal"
$ws.Range("A11").Value = "Religion Related"
$ws.Range("B11").Value = "'51.43%"
$ws.Range("B11").Style = "Normal"
$ws.Range("C11").Value = "'35"
$ws.Range("C11").Style = "Normal"
$ws.Range("D11").Value = "'`$12,350,058"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "'16.06%"
$ws.Range("E11").Style = "Normal"
$ws.Range("F11").Value = "'-0.30%"
$ws.Range("F11").Style = "Normal"
$ws.Range("A12").Value = "Unclassified"
$ws.Range("B12").Value = "'62.15%"
$ws.Range("B12").Style = "Normal"
$ws.Range("C12").Value = "'687"
$ws.Range("C12").Style = "Normal"
$ws.Range("D12").Value = "'`$1,942,850,201"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "'11.76%"
$ws.Range("E12").Style = "Normal"
$ws.Range("F12").Value = "'-7.16%"
$ws.Range("F12").Style = "Normal"
$ws.Range("A13").Value = "Universities"
$ws.Range("B13").Value = "'55.36%"
$ws.Range("B13").Style = "Normal"
$ws.Range("C13").Value = "'56"
$ws.Range("C13").Style = "Normal"
$ws.Range("D13").Value = "'`$2,076,673,548"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "'3.72%"
$ws.Range("E13").Style = "Normal"
$ws.Range("F13").Value = "'-0.89%"
$ws.Range("F13").Style = "Normal"
$ws.Range("A14").Value = "Total"
$ws.Range("B14").Value = "'60.30%"
$ws.Range("B14").Style = "Normal"
$ws.Range("C14").Value = "'3,018"
$ws.Range("C14").Style = "Normal"
$ws.Range("D14").Value = "'`$10,165,747,853"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "'11.74%"
$ws.Range("E14").Style = "Normal"
$ws.Range("F14").Value = "'-5.38%"
$ws.Range("F14").Style = "Normal"

